$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "conc S.D. [mg/l]" and "conc  [mg/l]" columns (B and C) were
# swapped, including their header text and all data values.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $cVal
    $ws.Cells.Item($r, 3).Value2 = $bVal
}
